$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "36.918.68"
$ws.Cells.Item(2, 5).Value = "  -0.56%  "

$ws.Cells.Item(3, 4).Value = "2.048.77"
$ws.Cells.Item(3, 5).Value = "  -0.14%  "

$ws.Cells.Item(4, 5).Value = "  -0.14%  "

$ws.Cells.Item(5, 4).Value = "'246.07"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.28%  "

$ws.Cells.Item(6, 4).Value = "'0.655"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -2.15%  "

$ws.Cells.Item(7, 4).Value = "'57.78"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -3.45%  "

$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(9, 4).Value = "'0.371"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -4.33%  "

$ws.Cells.Item(10, 4).Value = "'0.0779"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.85%  "

$ws.Cells.Item(11, 5).Value = "  +2.10%  "

$ws.Cells.Item(12, 4).Value = "'15.21"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -5.63%  "

$ws.Cells.Item(13, 4).Value = "'0.872"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +4.64%  "

$ws.Cells.Item(14, 4).Value = "2.349.49"
$ws.Cells.Item(14, 5).Value = "  -0.05%  "

$ws.Cells.Item(15, 2).Value = "WrappedEther"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(15, 4).Value = "2.344.46"
$ws.Cells.Item(15, 5).Value = "  +14.36%  "

$ws.Cells.Item(16, 2).Value = "Polkadot"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(16, 4).Value = "'5.61"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -3.09%  "

$ws.Cells.Item(17, 4).Value = "'17.93"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -2.11%  "

$ws.Cells.Item(18, 4).Value = "36.857.92"
$ws.Cells.Item(18, 5).Value = "  -0.76%  "

$ws.Cells.Item(19, 4).Value = "'73.50"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -3.51%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0888"
$ws.Cells.Item(20, 5).Value = "  -2.10%  "

$ws.Cells.Item(21, 4).Value = "'5.39"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.65%  "

$ws.Cells.Item(22, 4).Value = "'235.77"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.12%  "

$ws.Cells.Item(23, 5).Value = "  -0.01%  "

$ws.Cells.Item(24, 5).Value = "  +1.17%  "

$ws.Cells.Item(25, 5).Value = "  +9.39%  "

$ws.Cells.Item(26, 4).Value = "'2.19"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -1.19%  "

$ws.Cells.Item(27, 4).Value = "'168.42"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.59%  "

$ws.Cells.Item(28, 4).Value = "'19.93"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -1.62%  "

$ws.Cells.Item(29, 4).Value = "'5.45"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +12.90%  "

$ws.Cells.Item(30, 5).Value = "  -2.27%  "

$ws.Cells.Item(31, 5).Value = "  -4.02%  "

$ws.Cells.Item(32, 4).Value = "'4.69"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.58%  "

$ws.Cells.Item(33, 4).Value = "'0.0612"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -3.11%  "

$ws.Cells.Item(34, 4).Value = "'2.34"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.38%  "

$ws.Cells.Item(35, 5).Value = "  -0.05%  "

$ws.Cells.Item(36, 5).Value = "  +4.34%  "

$ws.Cells.Item(37, 4).Value = "'0.0825"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -8.04%  "

$ws.Cells.Item(38, 4).Value = "'1.31"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -2.52%  "

$ws.Cells.Item(39, 4).Value = "'5.14"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.34%  "

$ws.Cells.Item(40, 5).Value = "  -4.81%  "

$ws.Cells.Item(41, 5).Value = "  -0.83%  "

$ws.Cells.Item(42, 5).Value = "  -0.08%  "

$ws.Cells.Item(43, 4).Value = "'0.0946"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -13.67%  "

$ws.Cells.Item(44, 4).Value = "'96.62"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.84%  "

$ws.Cells.Item(45, 4).Value = "'16.88"
$ws.Cells.Item(45, 4).Style = "Normal"

$ws.Cells.Item(46, 4).Value = "1.297.57"
$ws.Cells.Item(46, 5).Value = "  +0.12%  "

$ws.Cells.Item(47, 4).Value = "'2.35"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -6.61%  "

$ws.Cells.Item(48, 5).Value = "  -0.82%  "

$ws.Cells.Item(49, 4).Value = "'6.73"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -1.90%  "

$ws.Cells.Item(50, 4).Value = "2.234.68"
$ws.Cells.Item(50, 5).Value = "  -0.37%  "

$ws.Cells.Item(51, 4).Value = "'44.88"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +1.23%  "
